$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: "Январь"/"Февраль"/"Май" columns collapse into a single
# "Апрель" column; the old Февраль/Май header cells become blank (keeping
# their style) and "Средняя температура" shifts left into column F.
$ws.Range("C1").Value = "Апрель"
$ws.Range("D1").ClearContents()
$ws.Range("E1").ClearContents()

# Row 2 data: city changes from Novosibirsk to Moscow, and the two sample
# values change; the average formula now spans C2:D2:E2 (incl. itself).
$ws.Range("B2").Value = "Москва"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("E2").Formula = "=AVERAGE(C2:D2:E2)"

# Column widths shrink to fit the new, shorter content.
$ws.Columns.Item(2).ColumnWidth = 6.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 7
$ws.Range("D1:E1").ColumnWidth = 1.1666666666666667
